$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.537.93'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.339.27'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.76'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '655.67'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.38'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -10.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.421'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -10.69%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.01'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.337.18'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.205'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -6.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.45'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '97.238.97'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.05'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000253'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -8.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.966.76'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.56'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.335.02'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.540'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +22.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.80'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.63'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '497.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -7.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.30'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -8.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000199'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -9.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.46'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.50'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -7.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.03'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -6.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.510.29'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.146'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.88'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.189'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.48'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +12.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.548'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.32'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.57'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.44'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.02%  '
$ws.Range("E41").Value = '  -7.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '506.08'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.57'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.70'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.835'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.63'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0412'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.48'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.64'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.36'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.14'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -11.23%  '
